# LMS-2523 Update BaSynthec Validation
# Update the strain identifier used in the OD600 bad-data example so it
# carries the "JJS-" prefix (MGP1000 -> JJS-MGP1000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-data")

$ws.Range("A3").Value = "JJS-MGP1000"
